# Applies the commit "Improved accuracy of stimulus presentation time-logging"
# Renames worksheets and updates the stimulus file-name values (column B)
# on each "task order" sheet, replacing stale timestamp-based filenames
# with newly generated ones.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new timestamp suffixes) ---------------------------
$wb.Worksheets.Item(1).Name = "GNG_TO-16512555502636547"
$wb.Worksheets.Item(2).Name = "NB_TO-1651255552429648"
$wb.Worksheets.Item(3).Name = "RS_TO-1651255552431649"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555525426533"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555526206465"

# --- Sheet 1: GNG_TO --------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555502226522.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555502456486.csv"
$ws1.Range("B4").Value = "go_stims-16512555502466486.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555502616491.csv"

# --- Sheet 2: NB_TO ----------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16512555507876475.csv"
$ws2.Range("B3").Value = "TB-16512555514716513.csv"
$ws2.Range("B4").Value = "ZB-match_5-16512555502856488.csv"
$ws2.Range("B5").Value = "ZB-match_8-16512555504056492.csv"
$ws2.Range("B6").Value = "TB-1651255552408648.csv"
$ws2.Range("B7").Value = "OB-16512555513486474.csv"
$ws2.Range("B8").Value = "OB-1651255550696649.csv"
$ws2.Range("B9").Value = "ZB-match_9-1651255550667648.csv"
$ws2.Range("B10").Value = "TB-16512555519926486.csv"

# --- Sheet 3: RS_TO -----------------------------------------------------------
# (only the sheet name changed; cell values are unchanged)

# --- Sheet 4: TOL_TO ----------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651255552461647.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555524387236.csv"
$ws4.Range("B4").Value = "MM_stims-16512555524776485.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555524636507.csv"
$ws4.Range("B6").Value = "MM_stims-16512555524936497.csv"
$ws4.Range("B7").Value = "ZM_stims-165125555247865.csv"

# --- Sheet 5: vSAT_TO ---------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16512555525476484.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555526056473.csv"
$ws5.Range("B4").Value = "vSAT_stims-1651255552588648.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555525726469.csv"
